# Session 1 notes update: append two new bullet paragraphs to the
# "Content Placeholder 2" shape on slide 2, after the existing
# "... or some other software." paragraph:
#   1) a normal bulleted paragraph introducing the GitHub repository
#   2) a no-bullet paragraph containing a hyperlinked URL

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

$url = "https://github.com/gitMarcH/Chanco_ST6103"

# Remember the length before we start inserting.
$lenBefore = $tr.Length

# New paragraph 1 - plain bulleted text at the same outline level (1)
# as the rest of the body text.
$tr.InsertAfter("`r" + "GitHub repository - will contain all course materials by the end of the week:") | Out-Null
$lenAfterPara1 = $tr.Length

# New paragraph 2 - the URL, with a trailing sentinel character so the
# range we manipulate below never sits at the very end of the text
# frame (doing so avoids the host materialising a stray <a:endParaRPr/>
# on the paragraph).
$tr.InsertAfter("`r" + $url + "X") | Out-Null
$lenAfterPara2 = $tr.Length

$urlStart = $lenAfterPara1 + 2
$urlLen   = $lenAfterPara2 - $lenAfterPara1 - 2

$urlRange = $tr.Characters($urlStart, $urlLen)
$urlRange.ActionSettings.Item(1).Hyperlink.Address = $url

# Drop the sentinel character again.
$tr.Characters($lenAfterPara2, 1).Text = ""

# Turn off the bullet for the URL-only paragraph.
$tr.Characters($urlStart, $urlLen).ParagraphFormat.Bullet.Visible = 0
